$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.357.63'
$ws.Range('E2').Value = '  +2.79%  '
$ws.Range('D3').Value = '3.812.96'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.73'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.84'
$ws.Range('E6').Value = '  +1.34%  '
$ws.Range('D7').Value = '3.809.36'
$ws.Range('E7').Value = '  +1.31%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.53'
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('E13').Value = '  -3.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.99'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').Value = '4.457.95'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '3.811.55'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '69.303.53'
$ws.Range('E17').Value = '  +2.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.30'
$ws.Range('E18').Value = '  -2.70%  '
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.16'
$ws.Range('E21').Value = '  +5.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '473.80'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.711'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.14'
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000149'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.25'
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.30'
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('D30').Value = '3.963.00'
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('E31').Value = '  -2.74%  '
$ws.Range('E32').Value = '  -2.43%  '
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.43'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.42'
$ws.Range('E35').Value = '  +3.13%  '
$ws.Range('D37').Value = '3.768.56'
$ws.Range('E38').Value = '  -2.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.60'
$ws.Range('E39').Value = '  -5.98%  '
$ws.Range('E40').Value = '  +2.13%  '
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.90'
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('E44').Value = '  -0.70%  '
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '44.23'
$ws.Range('E47').Value = '  +13.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.65'
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.43'
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '405.03'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '145.87'
$ws.Range('E51').Value = '  +2.90%  '
